# Update of correct answers in comprehension_questions.xlsx
# - CorrectResponse (column C) values corrected on several rows
# - the matching InfoRepeated/key text (column B / E) updated to the new key
# - the "kevesebb mint haromnegyedeben" explanatory texts simplified/reworded

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: "Mennyi időt vesz igénybe..." correct response b -> f
$ws.Range("C2").Value = "f"

# Row 3: "Megbizonyosodott róla..." answer keys swapped S/D -> F/J, correct response d -> j
$ws.Range("B3").Value = "F: Nem, nincs lehetőségem nyugodt körülmények között elvégezni a feladatot.    J: Igen, körülményeim nyugodtak, az internet kapcsolat stabil, a feladatokra tudok szánni 60 percet."
$ws.Range("C3").Value = "j"

# Row 6: correct response j -> k, repeated-info text reworded
$ws.Range("C6").Value = "k"
$ws.Range("E6").Value = "Kevesebb, mint 75%-ban tudott a képekre figyelni."

# Row 7: correct response j -> k, repeated-info text reworded (mixed formatting,
# second part keeps the original font/charset, matching the source edit)
$ws.Range("C7").Value = "k"
$ws.Range("E7").Value = "Kevesebb, mint 75%-ban tudott a képek helyszínére figyelni."
$ws.Range("E7").Font.Name = "Arial"
$ws.Range("E7").Font.Size = 10
$run2 = $ws.Range("E7").Characters(33, 28)
$run2.Font.Name = "Arial"
$run2.Font.Size = 10
